$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number (e.g. "217.81") would be
# auto-coerced to a numeric type by Excels normal ".Value" input parsing. The
# source data stores these as plain text (note sibling rows using "26.328.73"-style
# thousand-dotted text that safely stays text). Pre-marking the cell as Text via
# NumberFormat "@" before writing the value keeps it text, matching the original data.

$ws.Range("D2").Value = "26.328.73"
$ws.Range("E2").Value = "  -5.47%  "
$ws.Range("D3").Value = "1.671.63"
$ws.Range("E3").Value = "  -3.59%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.81"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5094"
$ws.Range("E6").Value = "  -11.31%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06365"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("E10").Value = "  -6.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07374"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "1.672.06"
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.558"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5828"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").Value = "1.899.12"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008524"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.88"
$ws.Range("E17").Value = "  -12.65%  "
$ws.Range("D18").Value = "26.395.43"
$ws.Range("E18").Value = "  -5.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.952"
$ws.Range("E19").Value = "  -6.49%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -3.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.32"
$ws.Range("E22").Value = "  -7.45%  "
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.42"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.684"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1180"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05906"
$ws.Range("E29").Value = "  -4.13%  "
$ws.Range("E30").Value = "  -8.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.321"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.517"
$ws.Range("E33").Value = "  -5.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.643"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("E36").Value = "  -6.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.355"
$ws.Range("E37").Value = "  -2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.651"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.025"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "1.077.13"
$ws.Range("E41").Value = "  -3.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8671"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.74"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "1.821.99"
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("E46").Value = "  +3.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.96"
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.066"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05190"
$ws.Range("E51").Value = "  -3.44%  "
